# Re-order the header columns in row 2 of the SubSector staging template.
# New order: SubSector_ID, BusinessKey, SectorBusinessKey, Code, LongName, ShortName, TextDescription

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "SubSector_ID"
$ws.Range("B2").Value = "BusinessKey"
$ws.Range("C2").Value = "SectorBusinessKey"
$ws.Range("D2").Value = "Code"
$ws.Range("E2").Value = "LongName"
$ws.Range("F2").Value = "ShortName"
$ws.Range("G2").Value = "TextDescription"
